# "origen en centro foto" - move field_x / field_y origin to the photo center.
# field_x (column K) and field_y (column L) are pixel coordinates measured
# from a corner of a 4656 x 3520 px photo; this recenters them on the photo
# center:
#   field_x_new = field_x_old - (width  / 2)
#   field_y_new = (height / 2) - field_y_old   (Y axis flip, image -> math/plot convention)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$width  = 4656
$height = 3520
$halfW  = $width  / 2
$halfH  = $height / 2

$lastRow = $ws.Cells.Item($ws.Rows.Count, 11).End(-4162).Row  # xlUp = -4162, col 11 = K

for ($row = 2; $row -le $lastRow; $row++) {
    $kCell = $ws.Cells.Item($row, 11)  # column K = field_x
    $lCell = $ws.Cells.Item($row, 12)  # column L = field_y

    $oldK = $kCell.Value2
    $oldL = $lCell.Value2

    if ($oldK -ne $null) {
        $kCell.Value2 = $oldK - $halfW
    }
    if ($oldL -ne $null) {
        $lCell.Value2 = $halfH - $oldL
    }
}
